$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.179.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.555.17'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.76%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '616.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.629'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.13%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.215'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.654'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.98'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000309'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.55'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.121.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '635.35'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +10.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '12.98'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.74%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.235.66'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.541.89'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.80%  '
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.57'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.74'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '103.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.53%  '
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.04'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.26%  '
$ws.Range('E27').Value = '  -2.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.07'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.32'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '63.89'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.70'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +19.14%  '
$ws.Range('E35').Value = '  -2.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '532.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.43%  '
$ws.Range('E37').Value = '  -3.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.32'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0782'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.533.15'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.25%  '
$ws.Range('E42').Value = '  +4.51%  '
$ws.Range('E43').Value = '  +1.53%  '
$ws.Range('E44').Value = '  +3.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.95'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.144'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.59%  '
$ws.Range('E47').Value = '  -5.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.22'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.33%  '
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.42'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '134.25'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.14%  '
